$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.408.60"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.478.23"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'489.25"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "'150.83"
$ws.Range("E6").Value = "  +7.11%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").Value = "2.486.61"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "'5.72"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("D11").Value = "'0.0981"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "2.908.59"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "56.660.22"
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "'21.11"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").Value = "2.498.41"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "'4.53"
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("D20").Value = "'10.24"
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("D21").Value = "'319.71"
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").Value = "'58.20"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "'0.160"
$ws.Range("E27").Value = "  -5.33%  "
$ws.Range("D28").Value = "2.597.37"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").Value = "'7.53"
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("D30").Value = "0.0₃0797"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "'151.36"
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "'18.24"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").Value = "'5.23"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").Value = "'0.865"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").Value = "'34.11"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("D41").Value = "'3.50"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").Value = "'0.0559"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'0.612"
$ws.Range("E43").Value = "  -0.34%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("D45").Value = "'265.58"
$ws.Range("E45").Value = "  +3.85%  "
$ws.Range("D46").Value = "'4.80"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").Value = "'0.0924"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("D48").Value = "'10.21"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("D49").Value = "'0.0228"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").Value = "'17.61"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").Value = "1.875.28"
$ws.Range("E51").Value = "  -6.38%  "
